# "automatically generating item via script"
# Adds a "Weight" and "Display Name" column to both BaseTable and TopTable
# on the "StuffToImport" sheet, repositions TopTable from E1:G4 to I1:M4,
# and removes the stray Name / LeatherBase_MetalHand rows (A7:A8).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("StuffToImport")

$baseTable = $ws.ListObjects.Item("BaseTable")
$topTable = $ws.ListObjects.Item("TopTable")

# Wipe everything in the working area - tables will be rebuilt from scratch
# in their new, final positions.
$ws.Range("A1:M8").ClearContents()

# ---- Base/Durability/Speed columns stay put (A:D) ----
$ws.Range("A1").Value = "Base"
$ws.Range("C1").Value = "Durability"
$ws.Range("D1").Value = "Speed"

$ws.Range("A2").Value = "WoodenBase"
$ws.Range("C2").Value = 10
$ws.Range("D2").Value = 5

$ws.Range("A3").Value = "MetalBase"
$ws.Range("C3").Value = 75
$ws.Range("D3").Value = 7

$ws.Range("A4").Value = "LeatherBase"
$ws.Range("C4").Value = 25
$ws.Range("D4").Value = 15

# ---- Top/Durability/Speed columns stay put (I:L) ----
$ws.Range("I1").Value = "Top"
$ws.Range("K1").Value = "Durability"
$ws.Range("L1").Value = "Speed"

$ws.Range("I2").Value = "WoodenHook"
$ws.Range("K2").Value = 5
$ws.Range("L2").Value = 5

$ws.Range("I3").Value = "MetalHook"
$ws.Range("K3").Value = 10
$ws.Range("L3").Value = 7

$ws.Range("I4").Value = "MetalHand"
$ws.Range("K4").Value = 15
$ws.Range("L4").Value = 8

# ---- Weight column added first (E for Base, M for Top) ----
$ws.Range("E1").Value = "Weight"
$ws.Range("E4").Value = 1

# Weight values that are text ("0.7", "1.2", "0.2", "0.5") rather than
# numbers - force a Text format so the assignment isn't auto-coerced to a
# number, then drop the formatting back to Normal so the cell keeps the
# string value without leaving a custom number format behind.
$textWeightCells = @(
    @{ addr = "E2"; value = "0.7" },
    @{ addr = "E3"; value = "1.2" },
    @{ addr = "M1"; value = "Weight" },
    @{ addr = "M2"; value = "0.2" },
    @{ addr = "M3"; value = "0.5" },
    @{ addr = "M4"; value = "0.7" }
)
foreach ($entry in $textWeightCells) {
    $ws.Range($entry.addr).NumberFormat = "@"
}
foreach ($entry in $textWeightCells) {
    $ws.Range($entry.addr).Value = $entry.value
}
foreach ($entry in $textWeightCells) {
    $ws.Range($entry.addr).Style = "Normal"
}

# ---- Display Name column added afterwards (B for Base, J for Top) ----
$ws.Range("B1").Value = "Display Name"
$ws.Range("B2").Value = "Wooden Base"
$ws.Range("B3").Value = "Metal Base"
$ws.Range("B4").Value = "Leather Base"

$ws.Range("J1").Value = "Display Name"
$ws.Range("J2").Value = "Wooden Hook"
$ws.Range("J3").Value = "Metal Hook"
$ws.Range("J4").Value = "Metal Hand"

# Resize the two tables onto their new ranges/columns.
$baseTable.Resize($ws.Range("A1:E4"))
$topTable.Resize($ws.Range("I1:M4"))

# Select J1, matching the saved selection in the target workbook.
$ws.Range("J1").Select()
